# Milestone 2 SCRUM content update:
# Reassign the "Login Page" user story (row 5) to Caleb Ljunggren instead
# of Ali Cooper, and leave the sheet's selection on that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = "Caleb Ljunggren"
$ws.Range("G5").Select()
